# "change name material group": rename the "Aluminum" material group value
# to "Aluminium" everywhere it is used (materialGroup / materialGroupChanged
# columns), and leave the BTMI002 sheet with cell H15 selected (matches the
# saved workbook view state).

$wb = $excel.ActiveWorkbook

$sheetsWithJ2 = @(
    "FPA001",
    "FPA002-003-005-007",
    "FPA004-006-010",
    "FPA008-009",
    "BTMI002",
    "BTMI003"
)

foreach ($name in $sheetsWithJ2) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("J2").Value = "Aluminium"
}

# BTMI015 has two "materialGroup" style columns (J = materialGroup,
# S = materialGroupChanged); on this sheet it's row 3 that was "Aluminum".
$wsBTMI015 = $wb.Worksheets.Item("BTMI015")
$wsBTMI015.Range("J3").Value = "Aluminium"
$wsBTMI015.Range("S3").Value = "Aluminium"

# Leave the BTMI002 sheet (the active tab) with H15 selected.
$wsBTMI002 = $wb.Worksheets.Item("BTMI002")
$wsBTMI002.Activate()
$wsBTMI002.Range("H15").Select()
